$d = $word.ActiveDocument
$failures = @()

function Replace-ExactText($findText, $newText) {
    $rng = $d.Content
    $found = $rng.Find.Execute($findText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if ($found) {
        $rng.Text = $newText
    }
    return $found
}

$ok0 = Replace-ExactText "Senior Data Engineer with over 20 years of enterprise experience, specializing in data analysis, system architecture, and technical problem-solving with proficiency in SQL, Python, Java, and C++. Adept at building scalable data infrastructure and optimizing workflows through data engineering and business intelligence practices. Experienced in leveraging relational databases and data science techniques to drive actionable insights and process improvements. Passionate about translating business goals into technical solutions, with a strong focus on operational reviews and product activation systems. Eager to contribute to YouTube’s mission of empowering creators and partners within Google’s innovative ecosystem by driving strategic initiatives for the partner ecosystem." "Senior Data Engineer with over 20 years of enterprise experience, specializing in data analysis, system architecture, and technical problem-solving with proficiency in SQL, Python, Java, and C++. Adept at building scalable data pipelines and infrastructure solutions, translating business goals into actionable technical systems, and driving process improvement through data-driven insights. Experienced in relational databases, data engineering, and business intelligence, with a proven track record of optimizing workflows and operational reviews. Passionate about contributing to YouTube’s mission of empowering creators and partners through innovative incentives systems and robust infrastructure development at Google. My expertise aligns with driving product activation at scale and collaborating with cross-functional teams to deliver impactful solutions. I am eager to leverage my skills to support YouTube’s global strategic initiatives."
if (-not $ok0) { $failures += "pair 0" }

$ok1 = Replace-ExactText "Architected automated ETL pipelines using Python and Pandas to ingest telemetry data from 6,000+ endpoints, enhancing data analysis and operational efficiency." "Architected automated ETL pipelines using Python and Pandas to ingest P95 telemetry from 6,000+ endpoints, enhancing data analysis and operational efficiency."
if (-not $ok1) { $failures += "pair 1" }

$ok2 = Replace-ExactText "Developed machine learning forecasting models with Prophet and scikit-learn to predict infrastructure bottlenecks six months ahead, supporting strategic planning." "Developed ML forecasting models with Prophet and scikit-learn to predict infrastructure bottlenecks 6 months ahead, improving provisioning accuracy."
if (-not $ok2) { $failures += "pair 2" }

$ok3 = Replace-ExactText "Designed optimized Oracle schemas for historical data retention, enabling accurate seasonal risk forecasting and system architecture improvements." "Designed optimized Oracle schemas for historical data retention, enabling seasonal risk forecasting and robust system architecture."
if (-not $ok3) { $failures += "pair 3" }

$ok4 = Replace-ExactText "Utilized SQL for data mining to identify underutilized infrastructure, driving hardware consolidation and significant cost savings." "Identified underutilized infrastructure through data mining, driving hardware consolidation and significant cost savings."
if (-not $ok4) { $failures += "pair 4" }

$ok5 = Replace-ExactText "Automated reporting workflows with Python scripts, streamlining process improvement and operational reviews." "Automated reporting workflows with Python scripts, streamlining process improvement and workflow optimization."
if (-not $ok5) { $failures += "pair 5" }

$ok6 = Replace-ExactText "Managed Dynatrace AppMon/Synthetics for critical systems, ensuring robust system architecture and performance monitoring." "Managed Dynatrace AppMon/Synthetics for Brand.com, focusing on performance data analysis for critical systems."
if (-not $ok6) { $failures += "pair 6" }

$ok7 = Replace-ExactText "Led 'FAST' project to data-mine real-user performance metrics, providing recommendations for system optimization and process improvement." "Led 'FAST' project to data-mine real-user performance metrics, providing optimization recommendations for system efficiency."
if (-not $ok7) { $failures += "pair 7" }

$ok8 = Replace-ExactText "Upgraded Dynatrace from 6.5 to 7.0, implemented TLS1.2 security, and supported cloud migration to AWS for enhanced infrastructure." "Upgraded DynaTrace (6.5 to 7.0) and supported cloud migration to AWS, enhancing system architecture and security (TLS1.2)."
if (-not $ok8) { $failures += "pair 8" }

$ok9 = Replace-ExactText "Developed dashboards for end-to-end functionality, delivering before/after metrics to support operational reviews." "Provided end-to-end monitoring and dashboarding, identifying bottlenecks for process improvement."
if (-not $ok9) { $failures += "pair 9" }

$ok10 = Replace-ExactText "Analyzed system bottlenecks using data analysis techniques and suggested performance enhancements for workflow optimization." "Analyzed performance issues and suggested actionable improvements to technical infrastructure."
if (-not $ok10) { $failures += "pair 10" }

$ok11 = Replace-ExactText "Integrated Performance Center with Dynatrace for comprehensive monitoring and data-driven insights." "Integrated Performance Center with DynaTrace for comprehensive monitoring solutions."
if (-not $ok11) { $failures += "pair 11" }

$ok12 = Replace-ExactText "Senior Consultant / SME for CA APM at CA Technologies (various clients) & Enterprise Iron (TIAA-CREF)" "SME for CA APM (Senior Consultant) at CA Technologies/ TIAA-CREF"
if (-not $ok12) { $failures += "pair 12" }

$ok13 = Replace-ExactText "Designed custom Management Modules, dashboards, and alerts using Perl/Ksh scripts for data extraction and operational reviews." "Designed custom Management Modules, dashboards, and alerts using Perl/Ksh scripts for data extraction and reporting."
if (-not $ok13) { $failures += "pair 13" }

$ok14 = Replace-ExactText "Provided sizing recommendations and Golden Images, enhancing infrastructure development and process improvement." "Provided sizing recommendations and bottleneck resolution for J2EE/.NET environments, enhancing technical infrastructure."
if (-not $ok14) { $failures += "pair 14" }

$ok15 = Replace-ExactText "Collaborated with IT teams to troubleshoot performance issues in J2EE/.NET environments, ensuring workflow optimization." "Collaborated with IT teams to troubleshoot performance issues, driving process improvement in complex systems."
if (-not $ok15) { $failures += "pair 15" }

$ok16 = Replace-ExactText "Trained client teams on APM solutions, fostering best practices in technical infrastructure and system monitoring." "Created Golden Images for agent rollouts, optimizing deployment workflows."
if (-not $ok16) { $failures += "pair 16" }

$ok17 = Replace-ExactText "Analyzed and resolved performance bottlenecks, delivering actionable insights for business intelligence." "Trained client teams on APM solutions, ensuring effective adoption and operational reviews."
if (-not $ok17) { $failures += "pair 17" }

$ok18 = Replace-ExactText "Analyzed J2EE telecom applications for load and break points, documenting key metrics like JDBC, threads, and memory for system architecture." "Analyzed J2EE telecom applications for load and breakpoints, focusing on data analysis of JDBC, threads, memory, CPU, and GC metrics."
if (-not $ok18) { $failures += "pair 18" }

$ok19 = Replace-ExactText "Installed JMX, Thread Dumps, and Wily Introscope to enhance performance monitoring and data analysis." "Installed JMX, Thread Dumps, and Wily Introscope to enhance system monitoring and architecture."
if (-not $ok19) { $failures += "pair 19" }

$ok20 = Replace-ExactText "Created automation scripts to streamline testing processes, supporting workflow optimization." "Created automation scripts to streamline performance testing processes."
if (-not $ok20) { $failures += "pair 20" }

$ok21 = Replace-ExactText "Identified resource bottlenecks through detailed data analysis, contributing to operational efficiency." "Documented key performance metrics to support technical problem-solving and system optimization."
if (-not $ok21) { $failures += "pair 21" }

$ok22 = Replace-ExactText " (2008-05 – 2012)" " (2008-05 – 2010)"
if (-not $ok22) { $failures += "pair 22" }

$ok23 = Replace-ExactText "Built CPPUNIT testing framework to automate conversion processes, enhancing workflow optimization." "Built CPPUNIT testing framework to automate conversion and ensure data integrity."
if (-not $ok23) { $failures += "pair 23" }

$ok24 = Replace-ExactText "Utilized data analysis to ensure performance metrics met business requirements during migration." "Enhanced technical infrastructure by refactoring database systems for scalability and performance."
if (-not $ok24) { $failures += "pair 24" }

$ok25 = Replace-ExactText "Collaborated with technical teams to design scalable infrastructure for high-performance systems." "Collaborated with teams to align migration with business goals and operational efficiency."
if (-not $ok25) { $failures += "pair 25" }

$ok26 = Replace-ExactText "Provided technical leadership in refactoring database systems for improved operational efficiency." "Provided data-driven insights to support decision-making during migration processes."
if (-not $ok26) { $failures += "pair 26" }

$ok27 = Replace-ExactText "Performed UML-based unit design for CICS/MQSeries/XML/DB2 systems, contributing to robust system architecture." "Conducted UML design and module development for a CICS/MQSeries/XML/DB2 system, focusing on system architecture."
if (-not $ok27) { $failures += "pair 27" }

$ok28 = Replace-ExactText "Collaborated on messaging architecture using VC++ and DB2 for operational efficiency." "Worked on messaging architecture using VC++ and DB2 for robust data handling."
if (-not $ok28) { $failures += "pair 28" }

$ok29 = Replace-ExactText "Developed high-availability multithreaded C++ interfaces using POSIX, sockets, and Marconi APIs for robust system architecture." "Developed high-availability multithreaded C++ interfaces using POSIX, sockets, and Marconi APIs for billing systems."
if (-not $ok29) { $failures += "pair 29" }

$ok30 = Replace-ExactText "Enhanced billing performance with C/C++/Pro*C/PL/SQL, achieving 75% memory reduction and 10x database performance." "Achieved 75% memory reduction and 20% throughput gain in billing processes through performance tuning with C/C++/Pro*C/PL/SQL."
if (-not $ok30) { $failures += "pair 30" }

$ok31 = Replace-ExactText "Automated system administration for WebLogic/WebSphere using Korn Shell scripts, supporting workflow optimization." "Automated system administration for WebLogic/WebSphere using Korn Shell scripts, driving workflow optimization."
if (-not $ok31) { $failures += "pair 31" }

$ok32 = Replace-ExactText "Designed interfaces with UML in Rational Rose, incorporating class/sequence diagrams for process improvement." "Designed interfaces with UML in Rational Rose, incorporating class/sequence diagrams and use cases."
if (-not $ok32) { $failures += "pair 32" }

$ok33 = Replace-ExactText "Troubleshot Enabler/CSM/EMS modules, ensuring operational reliability and efficiency." "Troubleshot Enabler/CSM/EMS modules, ensuring operational reliability and process improvement."
if (-not $ok33) { $failures += "pair 33" }

$ok34 = Replace-ExactText "Utilized SQL for database queries and reporting, driving actionable insights for business intelligence." "Implemented database fixes, improving performance by 10x through optimized sequences."
if (-not $ok34) { $failures += "pair 34" }

$ok35 = Replace-ExactText "Developed time/attendance interfaces using VB6/VC++ to support operational systems." "Developed time and attendance interfaces using VB6/VC++ to support business operations."
if (-not $ok35) { $failures += "pair 35" }

$ok36 = Replace-ExactText "Maintained and troubleshot Linux/Windows systems, performing various admin tasks for operational support." "Maintained and troubleshot Linux/Windows systems, performing various administrative tasks."
if (-not $ok36) { $failures += "pair 36" }

$ok37 = Replace-ExactText "Ensured system reliability through proactive monitoring and technical problem-solving." "Supported lab operations by ensuring system reliability and user support."
if (-not $ok37) { $failures += "pair 37" }

$ok38 = Replace-ExactText "SQL / Oracle, Python, Java, C++, ETL Design & Optimization, Data Warehousing, Pandas, scikit-learn, AWS, PySpark, Prophet / Time-Series Forecasting, Capacity Planning / Forecasting, GenAI / LLM Agents, Streamlit, PL/SQL, Perl, Ksh / Korn Shell Scripting, Dynatrace (AppMon + Synthetics), CA APM / Introscope, BMC TrueSight / TSCO, Oracle RAC, Multiprocessing, Docker, Git, Airflow, Hive/Hadoop, Linux/Unix, OCCI / OCI, WebLogic / WebSphere, VB6 / VC++" "SQL / Oracle, Python, C++, Java, Data Warehousing, ETL Design & Optimization, Pandas, scikit-learn, AWS, PySpark, Prophet / Time-Series Forecasting, Capacity Planning / Forecasting, GenAI / LLM Agents, Streamlit, PL/SQL, Perl, Ksh / Korn Shell Scripting, Dynatrace (AppMon + Synthetics), CA APM / Introscope, BMC TrueSight / TSCO, Oracle RAC, Multiprocessing, Docker, Git, Airflow, Hive/Hadoop, Linux/Unix, OCCI / OCI, WebLogic / WebSphere, VB6 / VC++"
if (-not $ok38) { $failures += "pair 38" }

$ok39 = Replace-ExactText "Optimized ETL and data access for enterprise-scale operations." "Optimized ETL and data access for enterprise scale."
if (-not $ok39) { $failures += "pair 39" }

$ok40 = Replace-ExactText "Created a modern agentic pipeline for banking-scale telemetry." "Built a modern agentic pipeline for banking-scale telemetry."
if (-not $ok40) { $failures += "pair 40" }

$ok41 = Replace-ExactText "Features comprehensive documentation, API guides, and performance benchmarks." "Repository includes code, documentation, and performance benchmarks."
if (-not $ok41) { $failures += "pair 41" }

$okSpecial = Replace-ExactText "^lUtilized Dynatrace for comprehensive analysis and actionable insights." ""
if (-not $okSpecial) { $failures += "special-br-removal" }

if ($failures.Count -gt 0) { Write-Host "FAILURES: $($failures -join ', ')" } else { Write-Host "All replacements applied successfully." }
